$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 655.6135653676964
$ws.Range("C3").Value = 551.4159515652706
$ws.Range("C4").Value = 551.7761419926437
$ws.Range("C5").Value = 628.0090048200738
$ws.Range("C6").Value = 627.8454998061012
$ws.Range("C7").Value = 644.8890145062114
$ws.Range("C8").Value = 704.1180029541127
$ws.Range("C9").Value = 669.4919776829566
$ws.Range("C10").Value = 655.8107634391241
$ws.Range("C11").Value = 688.5138766805632
$ws.Range("C12").Value = 691.7194150991568
$ws.Range("C13").Value = 679.5451382879403
$ws.Range("C14").Value = 686.526525763953
$ws.Range("C15").Value = 691.3562076339419
$ws.Range("C16").Value = 711.7736943531572
$ws.Range("C17").Value = 727.0002879594673
$ws.Range("C18").Value = 746.1512206659027
$ws.Range("C19").Value = 740.9403902338685
$ws.Range("C20").Value = 742.9816029688175
$ws.Range("C21").Value = 749.7743437626658
$ws.Range("C22").Value = 760.5751446438935
$ws.Range("C23").Value = 772.2627646198355
$ws.Range("C24").Value = 783.501310067658
$ws.Range("C25").Value = 790.5048296308337
$ws.Range("C26").Value = 797.9148360356403
$ws.Range("C27").Value = 804.8771327331394
$ws.Range("C28").Value = 814.9102111179029
$ws.Range("C29").Value = 821.0021946073454
$ws.Range("C30").Value = 825.097589907342
$ws.Range("C31").Value = 832.5428491231523
$ws.Range("C32").Value = 836.9803900389883
$ws.Range("C33").Value = 841.2458708373372
$ws.Range("C34").Value = 845.0759306302125
$ws.Range("C35").Value = 848.0956418046582
$ws.Range("C36").Value = 850.8196957273118
$ws.Range("C37").Value = 853.9952896683411
$ws.Range("C38").Value = 855.8230091790797
$ws.Range("C39").Value = 858.3959733194923
$ws.Range("C40").Value = 862.1891381531688
$ws.Range("C41").Value = 865.088800913072
$ws.Range("C42").Value = 865.6532136100958
$ws.Range("C43").Value = 867.8329600033674
$ws.Range("C44").Value = 869.3865399095336
$ws.Range("C45").Value = 870.7325170135283
$ws.Range("C46").Value = 871.2678564970488
$ws.Range("C47").Value = 872.3284442189172
$ws.Range("C48").Value = 873.8659577574026
$ws.Range("C49").Value = 874.6288808151314
$ws.Range("C50").Value = 876.3070160069386
$ws.Range("C51").Value = 878.567446600087
$ws.Range("C52").Value = 881.0688014435631
$ws.Range("C53").Value = 885.4207806928694
$ws.Range("C54").Value = 889.5972363174299
$ws.Range("C55").Value = 896.1835570983068
$ws.Range("C56").Value = 902.543097856201
$ws.Range("C57").Value = 905.4819096943056
$ws.Range("C58").Value = 906.8315150554396
$ws.Range("C59").Value = 908.4255157144617
$ws.Range("C60").Value = 914.2809693117871
$ws.Range("C61").Value = 925.8389506994203
$ws.Range("C62").Value = 927.4012527188402
$ws.Range("C63").Value = 928.9396604447164
$ws.Range("C64").Value = 930.1484398045457
$ws.Range("C65").Value = 931.3527854079018
